$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39 (sheet ALC)
$ws.Range("H39").Value = 489.7143
$ws.Range("I39").Value = 71.333336
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 214.000008
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 81.99999199999999
$ws.Range("N39").Value = -9592

# Row 64 (sheet ALC)
$ws.Range("H64").Value = 3874.1614
$ws.Range("I64").Value = 3712.5
$ws.Range("J64").Value = 4046.6
$ws.Range("K64").Value = 3712.5
$ws.Range("L64").Value = 4046.6
$ws.Range("M64").Value = -3464.5
$ws.Range("N64").Value = -4542.6

# Row 67 (sheet ALC)
$ws.Range("H67").Value = 3874.1614
$ws.Range("I67").Value = 3712.5
$ws.Range("J67").Value = 4046.6
$ws.Range("K67").Value = 3712.5
$ws.Range("L67").Value = 4046.6
$ws.Range("M67").Value = -2854.5
$ws.Range("N67").Value = -5762.6

# Row 125 (sheet ALC)
$ws.Range("H125").Value = 2041.2667
$ws.Range("J125").Value = 2041.2667
$ws.Range("L125").Value = 18371.4003
$ws.Range("N125").Value = -23291.4003

# Row 132 (sheet ALC)
$ws.Range("H132").Value = 2417.1853
$ws.Range("J132").Value = 2400
$ws.Range("L132").Value = 7200
$ws.Range("N132").Value = -12260

# Row 138 (sheet ALC)
$ws.Range("H138").Value = 3771.0344
$ws.Range("I138").Value = 1594.625
$ws.Range("J138").Value = 4600.143
$ws.Range("K138").Value = 4783.875
$ws.Range("L138").Value = 13800.429
$ws.Range("M138").Value = 356.125
$ws.Range("N138").Value = -24080.429

# Row 141 (sheet ALC)
$ws.Range("H141").Value = 1721.3334
$ws.Range("I141").Value = 1687.125
$ws.Range("K141").Value = 5061.375
$ws.Range("M141").Value = 118.625

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (sheet ARM)
$ws.Range("H74").Value = 3215.1628
$ws.Range("I74").Value = 4312.759
$ws.Range("J74").Value = 941.5714
$ws.Range("K74").Value = 4312.759
$ws.Range("L74").Value = 941.5714
$ws.Range("M74").Value = -3438.759
$ws.Range("N74").Value = -2689.5714

# Row 77 (sheet ARM)
$ws.Range("H77").Value = 3215.1628
$ws.Range("I77").Value = 4312.759
$ws.Range("J77").Value = 941.5714
$ws.Range("K77").Value = 21563.795
$ws.Range("L77").Value = 4707.857
$ws.Range("M77").Value = -17195.795
$ws.Range("N77").Value = -13443.857

# Row 110 (sheet ARM)
$ws.Range("H110").Value = 939
$ws.Range("I110").Value = 693.5333000000001
$ws.Range("J110").Value = 1859.5
$ws.Range("K110").Value = 693.5333000000001
$ws.Range("L110").Value = 1859.5
$ws.Range("M110").Value = 1351.4667
$ws.Range("N110").Value = -5949.5

# Row 132 (sheet ARM)
$ws.Range("H132").Value = 1321.0714
$ws.Range("I132").Value = 888.1389
$ws.Range("K132").Value = 2664.4167
$ws.Range("M132").Value = -134.4167000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (sheet BSM)
$ws.Range("H107").Value = 3470.6155
$ws.Range("I107").Value = 3976
$ws.Range("J107").Value = 3037.4285
$ws.Range("K107").Value = 3976
$ws.Range("L107").Value = 3037.4285
$ws.Range("M107").Value = -2056
$ws.Range("N107").Value = -6877.4285

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (sheet CRP)
$ws.Range("H16").Value = 2655.7222
$ws.Range("I16").Value = 1125.9
$ws.Range("J16").Value = 4568
$ws.Range("K16").Value = 1125.9
$ws.Range("L16").Value = 4568
$ws.Range("M16").Value = -838.9000000000001

# Row 99 (sheet CRP)
$ws.Range("H99").Value = 2389.3333
$ws.Range("I99").Value = 2058.182
$ws.Range("J99").Value = 3300
$ws.Range("K99").Value = 2058.182
$ws.Range("L99").Value = 3300
$ws.Range("M99").Value = -560.1819999999998
$ws.Range("N99").Value = -6296

# Row 107 (sheet CRP)
$ws.Range("H107").Value = 646.119
$ws.Range("I107").Value = 521.2083
$ws.Range("J107").Value = 812.6667
$ws.Range("K107").Value = 521.2083
$ws.Range("L107").Value = 812.6667
$ws.Range("M107").Value = 1398.7917
$ws.Range("N107").Value = -4652.6667

# Row 113 (sheet CRP)
$ws.Range("H113").Value = 2655.7222
$ws.Range("I113").Value = 1125.9
$ws.Range("J113").Value = 4568
$ws.Range("K113").Value = 1125.9
$ws.Range("L113").Value = 4568
$ws.Range("M113").Value = 1044.1

# Row 122 (sheet CRP)
$ws.Range("H122").Value = 1197.5834
$ws.Range("I122").Value = 965.46155
$ws.Range("J122").Value = 1471.909
$ws.Range("K122").Value = 2896.38465
$ws.Range("L122").Value = 4415.727000000001
$ws.Range("M122").Value = -446.38465
$ws.Range("N122").Value = -9315.727000000001

# Row 126 (sheet CRP)
$ws.Range("H126").Value = 2389.3333
$ws.Range("I126").Value = 2058.182
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 6174.545999999999
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -3704.545999999999
$ws.Range("N126").Value = -14840

# Row 132 (sheet CRP)
$ws.Range("H132").Value = 2513.25
$ws.Range("I132").Value = 1697.0526
$ws.Range("K132").Value = 5091.1578
$ws.Range("M132").Value = -2561.1578

# Row 134 (sheet CRP)
$ws.Range("H134").Value = 2887.6
$ws.Range("I134").Value = 1693.4667
$ws.Range("J134").Value = 6470
$ws.Range("K134").Value = 5080.4001
$ws.Range("L134").Value = 19410
$ws.Range("M134").Value = -2545.4001
$ws.Range("N134").Value = -24480.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 74 (sheet CUL)
$ws.Range("H74").Value = 1771
$ws.Range("I74").Value = 156.5
$ws.Range("K74").Value = 469.5
$ws.Range("M74").Value = 591.5

# Row 77 (sheet CUL)
$ws.Range("H77").Value = 1771
$ws.Range("I77").Value = 156.5
$ws.Range("K77").Value = 1408.5
$ws.Range("M77").Value = 3895.5

# Row 94 (sheet CUL)
$ws.Range("H94").Value = 2499.5
$ws.Range("I94").Value = 999
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 2997
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = -2321
$ws.Range("N94").Value = -13352

# Row 99 (sheet CUL)
$ws.Range("H99").Value = 1492.4286
$ws.Range("I99").Value = 808.5454999999999
$ws.Range("K99").Value = 2425.6365
$ws.Range("M99").Value = -179.6364999999996

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (sheet GSM)
$ws.Range("H102").Value = 1477.973
$ws.Range("I102").Value = 1515
$ws.Range("J102").Value = 1378
$ws.Range("K102").Value = 1515
$ws.Range("L102").Value = 1378
$ws.Range("M102").Value = 107
$ws.Range("N102").Value = -4622

# Row 107 (sheet GSM)
$ws.Range("H107").Value = 321.53845
$ws.Range("I107").Value = 331.75
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 331.75
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1588.25
$ws.Range("N107").Value = -4039

# Row 122 (sheet GSM)
$ws.Range("H122").Value = 3126453.5
$ws.Range("I122").Value = 4167938.2
$ws.Range("J122").Value = 1998.75
$ws.Range("K122").Value = 12503814.6
$ws.Range("L122").Value = 5996.25
$ws.Range("M122").Value = -12501364.6
$ws.Range("N122").Value = -10896.25

# Row 132 (sheet GSM)
$ws.Range("H132").Value = 2530.2666
$ws.Range("I132").Value = 2063.682
$ws.Range("J132").Value = 3813.375
$ws.Range("K132").Value = 6191.045999999999
$ws.Range("L132").Value = 11440.125
$ws.Range("M132").Value = -3661.045999999999
$ws.Range("N132").Value = -16500.125

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (sheet LTW)
$ws.Range("H7").Value = 2130.0454
$ws.Range("I7").Value = 2108.9443
$ws.Range("K7").Value = 2108.9443
$ws.Range("M7").Value = -1996.9443

# Row 126 (sheet LTW)
$ws.Range("H126").Value = 2130.0454
$ws.Range("I126").Value = 2108.9443
$ws.Range("K126").Value = 6326.8329
$ws.Range("M126").Value = -3856.8329

# Row 132 (sheet LTW)
$ws.Range("H132").Value = 6608.86
$ws.Range("I132").Value = 7781.737
$ws.Range("K132").Value = 23345.211
$ws.Range("M132").Value = -20815.211

# Row 136 (sheet LTW)
$ws.Range("H136").Value = 1712.325
$ws.Range("I136").Value = 1428.3334
$ws.Range("J136").Value = 2302.1538
$ws.Range("K136").Value = 4285.0002
$ws.Range("L136").Value = 6906.4614
$ws.Range("M136").Value = -1735.0002
$ws.Range("N136").Value = -12006.4614

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (sheet WVR)
$ws.Range("H107").Value = 530.65
$ws.Range("I107").Value = 559.05884
$ws.Range("J107").Value = 369.66666
$ws.Range("K107").Value = 1677.17652
$ws.Range("L107").Value = 1108.99998
$ws.Range("M107").Value = 242.82348
$ws.Range("N107").Value = -4948.999980000001

# Row 126 (sheet WVR)
$ws.Range("H126").Value = 1477.6562
$ws.Range("I126").Value = 1524.174
$ws.Range("J126").Value = 1358.7778
$ws.Range("K126").Value = 4572.522
$ws.Range("L126").Value = 4076.3334
$ws.Range("M126").Value = -2102.522
$ws.Range("N126").Value = -9016.3334

Write-Host "Edits applied"
